$wb = $excel.ActiveWorkbook

# --- Add the research note value to cells Q6:Q14 on the "Data" sheet ---
$wsData = $wb.Worksheets.Item("Data")
$noteValue = '"" '
for ($r = 6; $r -le 14; $r++) {
    $wsData.Range("Q" + $r).Value = $noteValue
}

# --- Switch the active sheet from "Research" to "Data" and update selection ---
$wsData.Activate()
$wsData.Range("R5").Select()

Write-Host "done"
